$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testcase")

# Update test case ids to include the LOGIN scope prefix
$ws.Range("A2").Value = "TC_LOGIN_0001"
$ws.Range("A3").Value = "TC_LOGIN_0002"
$ws.Range("A4").Value = "TC_LOGIN_0003"
$ws.Range("A5").Value = "TC_LOGIN_0004"

# Widen column A to fit the longer ids and move the active selection
$ws.Columns.Item(1).ColumnWidth = 16.43
$ws.Range("D22").Select()
